$d = $word.ActiveDocument

# --- 1. Merge "Inconsistent name chosen (eg. SW1 for a gate)." into a single run,
#        removing the proofErr (spell-check) markers that split it across runs.
$d.Content.Find.Execute(
    "Inconsistent name chosen (eg. SW1 for a gate).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Inconsistent name chosen (eg. SW1 for a gate).", 2) | Out-Null

# --- 2. Merge the "Look up name to find device type..." text into a single run,
#        removing the proofErr (grammar-check) markers around "no".
$d.Content.Find.Execute(
    "Look up name to find device type.  Check not SWITCH or CLOCK (which have no input).  For gates, check input number is in specified range.  Check only DTYPE uses DATA, CLOCK, SET, CLEAR.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Look up name to find device type.  Check not SWITCH or CLOCK (which have no input).  For gates, check input number is in specified range.  Check only DTYPE uses DATA, CLOCK, SET, CLEAR.", 2) | Out-Null

# --- 3. Append two new rows to the error-handling table.
$t = $d.Tables.Item(1)

$row1 = $t.Rows.Add()
$idx1 = $row1.Index
$t.Cell($idx1, 1).Range.Text = "Switch defined in 2 different states."
$t.Cell($idx1, 2).Range.Text = "Check if switch state has been defined before."

$row2 = $t.Rows.Add()
$idx2 = $row2.Index
$t.Cell($idx2, 1).Range.Text = "No monitor given."
$t.Cell($idx2, 2).Range.Text = "Check at least one monitor specified. If not, alert user."
